$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 347, shifting the existing rows 347-371
# down to 348-372 (dimension grows from A1:R371 to A1:R372).
$ws.Rows(347).Insert()

# Populate the newly inserted row 347 with the new weekly price record.
$ws.Range("A347").Value2 = 10
$ws.Range("B347").Value2 = "Vega Modelo de Temuco"
$ws.Range("C347").Value2 = "La Araucanía"
$ws.Range("D347").Value2 = 45223
$ws.Range("E347").Value2 = 9
$ws.Range("F347").Value2 = 100114007
$ws.Range("G347").Value2 = "Jengibre"
$ws.Range("H347").Value2 = "Sin especificar"
$ws.Range("I347").Value2 = "Primera"
$ws.Range("J347").Value2 = 50
$ws.Range("K347").Value2 = 32000
$ws.Range("L347").Value2 = 32000
$ws.Range("M347").Value2 = 32000
$ws.Range("N347").Value2 = "$/caja 13 kilos"
$ws.Range("O347").Value2 = "Perú"
$ws.Range("P347").Value2 = 2462
$ws.Range("Q347").Value2 = 13
$ws.Range("R347").Value2 = "Hortaliza"
